$d = $word.ActiveDocument

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$full = $d.Range($lastPara.Range.Start, $lastPara.Range.End)

$xmlFrag = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="52"/><w:szCs w:val="72"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="52"/><w:szCs w:val="72"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr><w:t>适合去看樱花</w:t></w:r></w:p>
<w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="52"/><w:szCs w:val="72"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="52"/><w:szCs w:val="72"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr><w:t>2023年3月9日星期四，多云</w:t></w:r></w:p>
<w:p><w:pPr><w:rPr><w:rFonts w:hint="default"/><w:sz w:val="52"/><w:szCs w:val="72"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="52"/><w:szCs w:val="72"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr><w:t>今天上午上了软件测试和开源项目课程</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$full.InsertXML($xmlFrag)
Write-Output "done"
